$d = $word.ActiveDocument

# The document contains a couple of unrelated (hidden, zero-content) text
# box shapes -- "Text Box 5" and "Text Box 6" -- that were reported to
# cause test failures on some systems. Remove them, along with the
# trailing run that held only a single space character right after them.

for ($i = $d.Shapes.Count; $i -ge 1; $i--) {
    $s = $d.Shapes.Item($i)
    if ($s.Name -eq "Text Box 5" -or $s.Name -eq "Text Box 6") {
        $s.Delete()
    }
}

# After removing the two shapes, the only remaining content in the main
# document flow is a run holding a single space character (which used to
# trail those shapes). Drop it too.
$r = $d.Range(0, 1)
if ($r.Text -eq " ") {
    $r.Delete()
}
